$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'55.394.44"
$ws.Range("E2").Value = "  -4.14%  "
$ws.Range("D3").Value = "'2.902.13"
$ws.Range("E3").Value = "  -4.50%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "'494.17"
$ws.Range("E5").Value = "  -4.07%  "
$ws.Range("D6").Value = "'132.04"
$ws.Range("E6").Value = "  -5.69%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.418"
$ws.Range("E8").Value = "  -5.46%  "
$ws.Range("D9").Value = "'7.14"
$ws.Range("E9").Value = "  -4.67%  "
$ws.Range("E10").Value = "  -6.69%  "
$ws.Range("D11").Value = "'0.348"
$ws.Range("E11").Value = "  -5.31%  "
$ws.Range("D12").Value = "'3.383.33"
$ws.Range("E12").Value = "  -4.94%  "
$ws.Range("D13").Value = "'0.125"
$ws.Range("E13").Value = "  -4.12%  "
$ws.Range("D14").Value = "'25.46"
$ws.Range("E14").Value = "  -4.51%  "
$ws.Range("D15").Value = "'0.0000156"
$ws.Range("E15").Value = "  -6.05%  "
$ws.Range("D16").Value = "'55.324.11"
$ws.Range("E16").Value = "  -4.34%  "
$ws.Range("D17").Value = "'5.97"
$ws.Range("E17").Value = "  -3.84%  "
$ws.Range("D18").Value = "'2.893.21"
$ws.Range("E18").Value = "  -4.65%  "
$ws.Range("D19").Value = "'12.41"
$ws.Range("E19").Value = "  -4.16%  "
$ws.Range("D20").Value = "'7.59"
$ws.Range("E20").Value = "  -5.23%  "
$ws.Range("D21").Value = "'312.50"
$ws.Range("E21").Value = "  -6.71%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "'0.480"
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'61.96"
$ws.Range("E24").Value = "  -4.46%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "'3.004.32"
$ws.Range("E25").Value = "  -5.15%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.158"
$ws.Range("E27").Value = "  -6.06%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "'0.0₃0838"
$ws.Range("E28").Value = "  -10.52%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'6.31"
$ws.Range("E29").Value = "  -7.51%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'6.93"
$ws.Range("E30").Value = "  -7.03%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.74"
$ws.Range("E31").Value = "  -4.35%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'19.55"
$ws.Range("E32").Value = "  -6.41%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.11"
$ws.Range("E33").Value = "  -8.61%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").Value = "'148.47"
$ws.Range("E34").Value = "  -4.81%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.35"
$ws.Range("E35").Value = "  -7.88%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'5.55"
$ws.Range("E36").Value = "  -5.32%  "
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").Value = "'24.24"
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.18"
$ws.Range("E38").Value = "  -8.17%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.0647"
$ws.Range("E39").Value = "  -5.83%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'36.18"
$ws.Range("E41").Value = "  -3.85%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'3.67"
$ws.Range("E42").Value = "  -5.38%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.630"
$ws.Range("E43").Value = "  -4.56%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.32"
$ws.Range("E44").Value = "  -7.85%  "
$ws.Range("D45").Value = "'2.076.40"
$ws.Range("E45").Value = "  -9.90%  "
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").Value = "'5.85"
$ws.Range("E46").Value = "  -3.22%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'0.909"
$ws.Range("E47").Value = "  -8.30%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0229"
$ws.Range("E48").Value = "  -4.58%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'18.51"
$ws.Range("E49").Value = "  -5.32%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.0836"
$ws.Range("E50").Value = "  -6.84%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'1.65"
$ws.Range("E51").Value = "  -9.60%  "
